$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.520.49"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "3.396.68"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'575.33"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").Value = "'138.28"
$ws.Range("E6").Value = "  +9.24%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "3.395.73"
$ws.Range("E8").Value = "  +3.86%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "'7.59"
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("E11").Value = "  +8.35%  "
$ws.Range("E12").Value = "  +6.58%  "
$ws.Range("D13").Value = "3.978.00"
$ws.Range("E13").Value = "  +3.83%  "
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("E15").Value = "  +8.33%  "
$ws.Range("D16").Value = "3.398.35"
$ws.Range("E16").Value = "  +3.47%  "
$ws.Range("D17").Value = "'25.45"
$ws.Range("E17").Value = "  +5.35%  "
$ws.Range("D18").Value = "61.593.09"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "'14.10"
$ws.Range("E19").Value = "  +6.72%  "
$ws.Range("D20").Value = "'5.91"
$ws.Range("E20").Value = "  +5.37%  "
$ws.Range("E21").Value = "  +4.31%  "
$ws.Range("D22").Value = "'390.11"
$ws.Range("E22").Value = "  +11.50%  "
$ws.Range("E23").Value = "  +4.06%  "
$ws.Range("D24").Value = "3.535.77"
$ws.Range("E24").Value = "  +3.88%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  +19.29%  "
$ws.Range("D27").Value = "'71.05"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("E28").Value = "  +14.66%  "
$ws.Range("D29").Value = "'7.75"
$ws.Range("E29").Value = "  +7.99%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'8.33"
$ws.Range("E31").Value = "  +7.27%  "
$ws.Range("D32").Value = "'0.161"
$ws.Range("E32").Value = "  +7.91%  "
$ws.Range("E33").Value = "  +3.41%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "3.428.16"
$ws.Range("E35").Value = "  +3.91%  "
$ws.Range("D36").Value = "'23.55"
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("E37").Value = "  +5.84%  "
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +5.72%  "
$ws.Range("D40").Value = "'161.83"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("E41").Value = "  +6.05%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.74"
$ws.Range("E42").Value = "  +12.51%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +8.74%  "
$ws.Range("E45").Value = "  +4.37%  "
$ws.Range("E46").Value = "  +5.26%  "
$ws.Range("D47").Value = "'41.21"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D48").Value = "'24.42"
$ws.Range("E48").Value = "  +9.90%  "
$ws.Range("D49").Value = "'6.99"
$ws.Range("E49").Value = "  +5.20%  "
$ws.Range("D50").Value = "'23.12"
$ws.Range("E50").Value = "  +9.01%  "
$ws.Range("D51").Value = "2.380.30"
$ws.Range("E51").Value = "  +10.87%  "